$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.546911358833313
$ws.Range("B1").Value = 2.627532958984375
$ws.Range("C1").Value = 3.133557081222534
$ws.Range("D1").Value = 3.725853204727173
$ws.Range("E1").Value = 1.374119162559509
